# Daily attendance processing - 2025-12-25 17:55:33
# Swap the first two comma-separated names/emails in the "Recorded By"
# column (G) for every data row on the active sheet, leaving any
# additional trailing entries (e.g. "system") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $first = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $first
            $newVal = $parts -join ", "
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
